# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G ("K" = strikeouts) is recomputed for every data row (rows 2-43)
# to replace the previous Strike# counts with the actual strikeout totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 5
    3  = 7
    4  = 5
    5  = 2
    6  = 5
    7  = 5
    8  = 7
    9  = 3
    10 = 3
    11 = 2
    12 = 4
    13 = 2
    14 = 5
    15 = 5
    16 = 4
    17 = 7
    18 = 3
    19 = 3
    20 = 3
    21 = 3
    22 = 3
    23 = 1
    24 = 0
    25 = 5
    26 = 1
    27 = 0
    28 = 0
    29 = 0
    30 = 3
    31 = 4
    32 = 2
    33 = 0
    34 = 1
    35 = 2
    36 = 1
    37 = 1
    38 = 0
    39 = 0
    40 = 4
    41 = 5
    42 = 3
    43 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
